# Auto-generated edit script applying the Moogle_Profits diff
# Updates currentAveragePrice / LevePrice / LeveProfit derived columns (H,I,J,K,L,M,N)
# across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets per the commit's refreshed market data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 1475
$ws.Range("I5").Value = 1572.3684
$ws.Range("K5").Value = 1572.3684
$ws.Range("M5").Value = -1457.3684
$ws.Range("H11").Value = 493.54544
$ws.Range("I11").Value = 493.54544
$ws.Range("K11").Value = 493.54544
$ws.Range("M11").Value = -353.54544
$ws.Range("H53").Value = 853.5625
$ws.Range("J53").Value = 1409.8334
$ws.Range("L53").Value = 1409.8334
$ws.Range("N53").Value = -2683.8334
$ws.Range("H69").Value = 17908.273
$ws.Range("I69").Value = 12995
$ws.Range("J69").Value = 18399.6
$ws.Range("K69").Value = 38985
$ws.Range("L69").Value = 55198.8
$ws.Range("M69").Value = -38111
$ws.Range("N69").Value = -56946.8
$ws.Range("H72").Value = 17908.273
$ws.Range("I72").Value = 12995
$ws.Range("J72").Value = 18399.6
$ws.Range("K72").Value = 116955
$ws.Range("L72").Value = 165596.4
$ws.Range("M72").Value = -112587
$ws.Range("N72").Value = -174332.4
$ws.Range("H86").Value = 4840.3
$ws.Range("I86").Value = 1466.1666
$ws.Range("J86").Value = 6286.357
$ws.Range("K86").Value = 1466.1666
$ws.Range("L86").Value = 6286.357
$ws.Range("M86").Value = -343.1666
$ws.Range("N86").Value = -8532.357
$ws.Range("H89").Value = 4840.3
$ws.Range("I89").Value = 1466.1666
$ws.Range("J89").Value = 6286.357
$ws.Range("K89").Value = 7330.833000000001
$ws.Range("L89").Value = 31431.785
$ws.Range("M89").Value = -1714.833000000001
$ws.Range("N89").Value = -42663.785
$ws.Range("H92").Value = 502.93332
$ws.Range("I92").Value = 568.5
$ws.Range("K92").Value = 568.5
$ws.Range("M92").Value = 679.5
$ws.Range("H132").Value = 3254.535
$ws.Range("I132").Value = 3235.5527
$ws.Range("K132").Value = 9706.658100000001
$ws.Range("M132").Value = -7176.658100000001
$ws.Range("H135").Value = 583.2727
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H137").Value = 3043.913
$ws.Range("I137").Value = 1459.125
$ws.Range("J137").Value = 3889.1333
$ws.Range("K137").Value = 4377.375
$ws.Range("L137").Value = 11667.3999
$ws.Range("M137").Value = -1827.375
$ws.Range("N137").Value = -16767.3999
$ws.Range("H138").Value = 3109.7576
$ws.Range("I138").Value = 2494.8572
$ws.Range("J138").Value = 6553.2
$ws.Range("K138").Value = 7484.571599999999
$ws.Range("L138").Value = 19659.6
$ws.Range("M138").Value = -2344.571599999999
$ws.Range("N138").Value = -29939.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5786.6143
$ws.Range("I32").Value = 2601.2153
$ws.Range("J32").Value = 47196.8
$ws.Range("K32").Value = 2601.2153
$ws.Range("L32").Value = 47196.8
$ws.Range("M32").Value = -2314.2153
$ws.Range("N32").Value = -47770.8
$ws.Range("H61").Value = 9971.964
$ws.Range("I61").Value = 7903.2104
$ws.Range("J61").Value = 14339.333
$ws.Range("K61").Value = 7903.2104
$ws.Range("L61").Value = 14339.333
$ws.Range("M61").Value = -7691.2104
$ws.Range("N61").Value = -14763.333
$ws.Range("H74").Value = 4080.516
$ws.Range("I74").Value = 2094.76
$ws.Range("K74").Value = 2094.76
$ws.Range("M74").Value = -1220.76
$ws.Range("H77").Value = 4080.516
$ws.Range("I77").Value = 2094.76
$ws.Range("K77").Value = 10473.8
$ws.Range("M77").Value = -6105.800000000001
$ws.Range("H132").Value = 3842.879
$ws.Range("I132").Value = 1816.8
$ws.Range("J132").Value = 6959.923
$ws.Range("K132").Value = 5450.4
$ws.Range("L132").Value = 20879.769
$ws.Range("M132").Value = -2920.4
$ws.Range("N132").Value = -25939.769
$ws.Range("H136").Value = 9971.964
$ws.Range("I136").Value = 7903.2104
$ws.Range("J136").Value = 14339.333
$ws.Range("K136").Value = 23709.6312
$ws.Range("L136").Value = 43017.999
$ws.Range("M136").Value = -21159.6312
$ws.Range("N136").Value = -48117.999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H21").Value = 52449.375
$ws.Range("J21").Value = 52449.375
$ws.Range("L21").Value = 52449.375
$ws.Range("N21").Value = -52921.375
$ws.Range("H54").Value = 32846
$ws.Range("J54").Value = 45498.25
$ws.Range("L54").Value = 45498.25
$ws.Range("N54").Value = -46466.25
$ws.Range("H86").Value = 2131.625
$ws.Range("I86").Value = 1719.8948
$ws.Range("K86").Value = 1719.8948
$ws.Range("M86").Value = -596.8948
$ws.Range("H89").Value = 2131.625
$ws.Range("I89").Value = 1719.8948
$ws.Range("K89").Value = 8599.474
$ws.Range("M89").Value = -2983.474
$ws.Range("H99").Value = 4322.6875
$ws.Range("I99").Value = 2955.4
$ws.Range("J99").Value = 6601.5
$ws.Range("K99").Value = 2955.4
$ws.Range("L99").Value = 6601.5
$ws.Range("M99").Value = -1457.4
$ws.Range("N99").Value = -9597.5
$ws.Range("H134").Value = 2824.3333
$ws.Range("I134").Value = 1619.6296
$ws.Range("J134").Value = 13666.667
$ws.Range("K134").Value = 4858.8888
$ws.Range("L134").Value = 41000.001
$ws.Range("M134").Value = -2323.8888
$ws.Range("N134").Value = -46070.001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8152.4165
$ws.Range("I31").Value = 4034.7
$ws.Range("J31").Value = 11093.643
$ws.Range("K31").Value = 4034.7
$ws.Range("L31").Value = 11093.643
$ws.Range("M31").Value = -3739.7
$ws.Range("N31").Value = -11683.643
$ws.Range("H34").Value = 8152.4165
$ws.Range("I34").Value = 4034.7
$ws.Range("J34").Value = 11093.643
$ws.Range("K34").Value = 4034.7
$ws.Range("L34").Value = 11093.643
$ws.Range("M34").Value = -3832.7
$ws.Range("N34").Value = -11497.643
$ws.Range("H38").Value = 17499.5
$ws.Range("I38").Value = 4999.5
$ws.Range("K38").Value = 4999.5
$ws.Range("M38").Value = -4622.5
$ws.Range("H46").Value = 17499.5
$ws.Range("I46").Value = 4999.5
$ws.Range("K46").Value = 4999.5
$ws.Range("M46").Value = -4788.5
$ws.Range("H132").Value = 3378.6316
$ws.Range("I132").Value = 2548.7407
$ws.Range("J132").Value = 5415.636
$ws.Range("K132").Value = 7646.222099999999
$ws.Range("L132").Value = 16246.908
$ws.Range("M132").Value = -5116.222099999999
$ws.Range("N132").Value = -21306.908
$ws.Range("H133").Value = 60000
$ws.Range("J133").Value = 60000
$ws.Range("L133").Value = 60000
$ws.Range("N133").Value = -65060
$ws.Range("H134").Value = 4674.4375
$ws.Range("I134").Value = 2789.7
$ws.Range("K134").Value = 8369.099999999999
$ws.Range("M134").Value = -5834.099999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H21").Value = 266
$ws.Range("I21").Value = 257.5
$ws.Range("J21").Value = 300
$ws.Range("K21").Value = 772.5
$ws.Range("L21").Value = 900
$ws.Range("M21").Value = -599.5
$ws.Range("N21").Value = -1246
$ws.Range("H113").Value = 1164.2727
$ws.Range("I113").Value = 401
$ws.Range("J113").Value = 1800.3334
$ws.Range("K113").Value = 1203
$ws.Range("L113").Value = 5401.0002
$ws.Range("M113").Value = 967
$ws.Range("N113").Value = -9741.0002
$ws.Range("H122").Value = 817.1818
$ws.Range("J122").Value = 798.5714
$ws.Range("L122").Value = 7187.1426
$ws.Range("N122").Value = -12087.1426

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 84000
$ws.Range("J123").Value = 84000
$ws.Range("L123").Value = 84000
$ws.Range("N123").Value = -88900
$ws.Range("H132").Value = 5738.7124
$ws.Range("I132").Value = 5667.7417
$ws.Range("K132").Value = 17003.2251
$ws.Range("M132").Value = -14473.2251

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4698.923
$ws.Range("I40").Value = 3208.1
$ws.Range("K40").Value = 3208.1
$ws.Range("M40").Value = -3072.1
$ws.Range("H55").Value = 1132.25
$ws.Range("I55").Value = 570.1111
$ws.Range("K55").Value = 570.1111
$ws.Range("M55").Value = -397.1111
$ws.Range("H82").Value = 1033.2632
$ws.Range("I82").Value = 638.5714
$ws.Range("K82").Value = 638.5714
$ws.Range("M82").Value = -277.5714
$ws.Range("H85").Value = 1033.2632
$ws.Range("I85").Value = 638.5714
$ws.Range("K85").Value = 638.5714
$ws.Range("M85").Value = 609.4286
$ws.Range("H100").Value = 4796.357
$ws.Range("I100").Value = 4579.778
$ws.Range("K100").Value = 4579.778
$ws.Range("M100").Value = -4038.778
$ws.Range("H132").Value = 4782.025
$ws.Range("I132").Value = 4281.706
$ws.Range("J132").Value = 7617.1665
$ws.Range("K132").Value = 12845.118
$ws.Range("L132").Value = 22851.4995
$ws.Range("M132").Value = -10315.118
$ws.Range("N132").Value = -27911.4995
$ws.Range("H136").Value = 5283.0303
$ws.Range("I136").Value = 4487.25
$ws.Range("J136").Value = 6032
$ws.Range("K136").Value = 13461.75
$ws.Range("L136").Value = 18096
$ws.Range("M136").Value = -10911.75
$ws.Range("N136").Value = -23196
$ws.Range("H139").Value = 62499.5
$ws.Range("I139").Value = 25000
$ws.Range("K139").Value = 25000
$ws.Range("M139").Value = -19860

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 16999.5
$ws.Range("J74").Value = 18333
$ws.Range("L74").Value = 18333
$ws.Range("N74").Value = -20205
$ws.Range("H77").Value = 16999.5
$ws.Range("J77").Value = 18333
$ws.Range("L77").Value = 54999
$ws.Range("N77").Value = -64359
$ws.Range("H102").Value = 199500
$ws.Range("J102").Value = 199500
$ws.Range("L102").Value = 199500
$ws.Range("N102").Value = -205990
$ws.Range("H132").Value = 2359.4285
$ws.Range("I132").Value = 1577.4
$ws.Range("J132").Value = 18000
$ws.Range("K132").Value = 4732.200000000001
$ws.Range("L132").Value = 54000
$ws.Range("M132").Value = -2202.200000000001
$ws.Range("N132").Value = -59060
